# This edit refreshes a scraped product export:
#  - the whole sheet gets a newer scrape "timestamp" (column O) for every
#    data row (the site was re-crawled later the same day)
#  - a handful of rows were re-scraped out of order, so their whole content
#    (every column) ends up swapped / rotated between a few row numbers
#
# Columns (1-based): A id, B title, C href, D quantity, E ratingAmount,
# F ratingValue, G brand, H price, I priceContext, J priceContextHiddenText,
# K priceContextPrice, L priceContextAmount, M udoCat, N productAriaLabel,
# O timestamp
#
# Only E (ratingAmount) and F (ratingValue) are genuine numbers in this
# sheet; every other column holds text (even when it looks numeric, e.g.
# price "16.95" or id "4014528"), so those columns must be forced to Text
# number format before the value is written back, otherwise Excel would
# "helpfully" convert strings like "1.70" into the number 1.7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 15
$numericCols = @(5, 6)   # E = ratingAmount, F = ratingValue

function Get-RowValues($row) {
    $vals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $vals += ,$ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($row, $c)
        if ($numericCols -contains $c) {
            $cell.NumberFormat = "General"
        } else {
            $cell.NumberFormat = "@"
        }
        $v = $vals[$c - 1]
        if ($null -eq $v) {
            $cell.Value = ""
        } else {
            $cell.Value = $v
        }
    }
}

function Swap-Rows($rowA, $rowB) {
    $a = Get-RowValues $rowA
    $b = Get-RowValues $rowB
    Set-RowValues $rowA $b
    Set-RowValues $rowB $a
}

# --- Rearrange the rows whose content was re-scraped out of order ---

# Simple pairwise swaps
Swap-Rows 8 9
Swap-Rows 18 19
Swap-Rows 30 31

# Rows 21-24 rotate by one: 21<-24, 22<-21, 23<-22, 24<-23
$r21 = Get-RowValues 21
$r22 = Get-RowValues 22
$r23 = Get-RowValues 23
$r24 = Get-RowValues 24
Set-RowValues 21 $r24
Set-RowValues 22 $r21
Set-RowValues 23 $r22
Set-RowValues 24 $r23

# Rows 54-56 rotate by one: 54<-56, 55<-54, 56<-55
$r54 = Get-RowValues 54
$r55 = Get-RowValues 55
$r56 = Get-RowValues 56
Set-RowValues 54 $r56
Set-RowValues 55 $r54
Set-RowValues 56 $r55

# --- Refresh the scrape timestamp on every data row ---
$newTimestamp = "2022-08-24 20:57:41"
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 86 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 15)
    $cell.NumberFormat = "@"
    $cell.Value = $newTimestamp
}
